$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for columns A,B,D,E,F,G,H,Q,R for rows 2..13
# before making any writes, since the update is a permutation of rows.
$cols = @(1,2,4,5,6,7,8,17,18)
$snapshot = @{}
for ($r = 2; $r -le 13; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowData
}

# Map: target row -> source row (where the new row content comes from)
$mapping = @{
    2  = 3
    3  = 10
    4  = 6
    5  = 12
    6  = 4
    7  = 11
    8  = 2
    9  = 5
    10 = 7
    11 = 8
    12 = 13
    13 = 9
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $src = $snapshot[$source]

    $ws.Cells.Item($target, 1).Value  = $src[1]
    $ws.Cells.Item($target, 2).Value  = $src[2]
    $ws.Cells.Item($target, 4).Value  = $src[4]
    $ws.Cells.Item($target, 5).Value  = $src[5]
    $ws.Cells.Item($target, 6).Value  = $src[6]
    $ws.Cells.Item($target, 7).Value  = $src[7]
    $ws.Cells.Item($target, 8).Value  = $src[8]
    $ws.Cells.Item($target, 17).Value = [Math]::Round($src[17])
    $ws.Cells.Item($target, 18).Value = [Math]::Round($src[18])
}
